# Update the error rows: clear the Data_Emissao (column K) timestamp and
# change the Mensagem_Erro (column L) text for rows 2-5.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 5; $row++) {
    $ws.Cells.Item($row, 11).ClearContents()              # K column - Data_Emissao
    $ws.Cells.Item($row, 12).Value = "Erro ao selecionar atividade"   # L column - Mensagem_Erro
}
